$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (new entry #7): copy date/percent formatted template cells first so
# the reused numFmt/style indices match existing ones (avoids creating new
# numFmt/style entries), then overwrite with the real values.
$ws.Range("B10").Copy($ws.Range("B13"))
$ws.Range("E10").Copy($ws.Range("E13"))

$ws.Range("A13").Value = 7
$ws.Range("B13").Value = 44572
$ws.Range("C13").Value = "RPA RLOGIC"
$ws.Range("D13").Value = "1. Generated the Accounting Statements for the three centers and shared to Rahaman san to verify"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Completed"

# Row 14 (continuation line): same percent-style template copy for E14.
$ws.Range("E8").Copy($ws.Range("E14"))

$ws.Range("D14").Value = "2. Uploading the pdf files task is work in progress for ESA "
$ws.Range("E14").Value = 0.1
$ws.Range("F14").Value = "WIP"

# Move the active selection from D21 to D14, as in the author's edit.
$ws.Range("D14").Select()
